$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "table: category" -> "table: Category"
$ws.Range("A7").Value = "table: Category"

# "table: level" -> "table: Level"
$ws.Range("E7").Value = "table: Level"

# Update selection to E7:F7 (active cell E7)
$null = $ws.Range("E7:F7").Select()
